$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename the group from "质控组" to "北京组" on both sheets
$ws1.Range("A2:A5").Value = "北京组"
$ws2.Range("A2").Value = "北京组"

# Updated metrics on Sheet2 row 2
$ws2.Range("G2").Value = 23
$ws2.Range("H2").Value = 26
$ws2.Range("I2").Value = 5.07

# Restore the cursor/selection position on each sheet
[void]$ws2.Range("A2").Select()
[void]$ws1.Range("A5").Select()
